$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh — Fri Sep  6 16:51:07 UTC 2024 run

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '53.953.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.259.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.57%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '487.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '126.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.93%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.520'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.260.08'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0923'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.82%  '
$ws.Range("E11").Value = '  -1.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.78'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.315'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.661.92'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '53.870.95'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000128'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.260.14'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.87%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.30%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '300.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.10'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.28%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.365'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.143'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '169.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0691'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.06'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.44'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.16'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.831'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.59'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.79'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.366'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.36'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.76'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0876'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.538'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '237.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0472'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0202'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.18%  '
